$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# ALC sheet - row 137
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1715700.8
$ws.Range("J137").Value = 7579.8823
$ws.Range("L137").Value = 22739.6469
$ws.Range("N137").Value = -27839.6469

# ---------------------------------------------------------------------------
# ARM sheet - row 123
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H123").Value = 40429
$ws.Range("J123").Value = 40429
$ws.Range("L123").Value = 40429
$ws.Range("N123").Value = -50229

# ---------------------------------------------------------------------------
# CUL sheet - rows 68, 71, 129
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H68").Value = 761487.2
$ws.Range("I68").Value = 756979.1
$ws.Range("J68").Value = 762784
$ws.Range("K68").Value = 2270937.3
$ws.Range("L68").Value = 2288352
$ws.Range("M68").Value = -2270126.3
$ws.Range("N68").Value = -2289974

$ws.Range("H71").Value = 761487.2
$ws.Range("I71").Value = 756979.1
$ws.Range("J71").Value = 762784
$ws.Range("K71").Value = 6812811.899999999
$ws.Range("L71").Value = 6865056
$ws.Range("M71").Value = -6808755.899999999
$ws.Range("N71").Value = -6873168

$ws.Range("H129").Value = 92385.55
$ws.Range("I129").Value = 334146.22
$ws.Range("J129").Value = 1725.2916
$ws.Range("K129").Value = 1002438.66
$ws.Range("L129").Value = 5175.8748
$ws.Range("M129").Value = -997438.6599999999
$ws.Range("N129").Value = -15175.8748

# ---------------------------------------------------------------------------
# LTW sheet - row 60
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H60").Value = 33296.668
$ws.Range("J60").Value = 33296.668
$ws.Range("L60").Value = 33296.668
$ws.Range("N60").Value = -34314.668

# ---------------------------------------------------------------------------
# WVR sheet - rows 119-141 (new H:N values)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H119").Value = 40869.7
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 40869.7
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 40869.7
$ws.Range("N119").Value = -50545.7

$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 0

$ws.Range("H121").Value = 30420
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 30420
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 30420
$ws.Range("N121").Value = -33914

$ws.Range("H122").Value = 1906624.2
$ws.Range("I122").Value = 2382964.2
$ws.Range("J122").Value = 1263.6666
$ws.Range("K122").Value = 7148892.600000001
$ws.Range("L122").Value = 3790.9998
$ws.Range("M122").Value = -7146442.600000001
$ws.Range("N122").Value = -8690.9998

$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 0

$ws.Range("H124").Value = 45429
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 45429
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 45429
$ws.Range("N124").Value = -55249

$ws.Range("H125").Value = 36633.332
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 36633.332
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 36633.332
$ws.Range("N125").Value = -46473.332

$ws.Range("H126").Value = 1090718
$ws.Range("I126").Value = 1177775.4
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 3533326.2
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -3530856.2
$ws.Range("N126").Value = -12440

$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 0

$ws.Range("H128").Value = 48000
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 48000
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 48000
$ws.Range("N128").Value = -57960

$ws.Range("H129").Value = 40429
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 40429
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 40429
$ws.Range("N129").Value = -50429

$ws.Range("H130").Value = 30428.5
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 30428.5
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 30428.5
$ws.Range("N130").Value = -40468.5

$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 0

$ws.Range("H132").Value = 1064724.5
$ws.Range("I132").Value = 1557713
$ws.Range("J132").Value = 2903.1538
$ws.Range("K132").Value = 4673139
$ws.Range("L132").Value = 8709.4614
$ws.Range("M132").Value = -4670609
$ws.Range("N132").Value = -13769.4614

$ws.Range("H133").Value = 50531.168
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 50531.168
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 50531.168
$ws.Range("N133").Value = -60651.168

$ws.Range("H135").Value = 38990.266
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 38990.266
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 38990.266
$ws.Range("N135").Value = -49130.266

$ws.Range("H136").Value = 1015905.2
$ws.Range("I136").Value = 1667738
$ws.Range("J136").Value = 1943.1111
$ws.Range("K136").Value = 5003214
$ws.Range("L136").Value = 5829.3333
$ws.Range("M136").Value = -5000664
$ws.Range("N136").Value = -10929.3333

$ws.Range("H137").Value = 47000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 47000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 47000
$ws.Range("N137").Value = -57200

$ws.Range("H138").Value = 46000
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 46000
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 46000
$ws.Range("N138").Value = -56280

$ws.Range("H139").Value = 59000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 59000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 59000
$ws.Range("N139").Value = -69280

$ws.Range("H140").Value = 41586.555
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 41586.555
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 41586.555
$ws.Range("N140").Value = -51946.555

$ws.Range("H141").Value = 39200
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 39200
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 39200
$ws.Range("N141").Value = -49560

Write-Output "Edits applied"
